# Fixed naive component forecaster bug - Presentation state 11.02.
#
# The naive QoQ error series (row-wise, one row per forecast-origin quarter)
# was missing its first forecast-horizon column: every row's error series
# needs a new value inserted at column B (the "1-quarter-ahead" bucket),
# with the rest of that row's values shifting one column to the right
# (columns C..K), matching the newly-added "matched to ifoCAST" horizon.
# The sheet's used range is fixed at column K, so for the long rows
# (2-10, which already reach column K) the old right-most value is
# pushed out; shorter rows simply grow by one column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value to insert at column B for each data row (2-20).
$newB = @{
  2  = "0.2433086034063205"
  3  = "-0.115952868393072"
  4  = "0.1459581181302581"
  5  = "-0.08819670345554087"
  6  = "0.3847923593882046"
  7  = "0.03766489642184559"
  8  = "0.1769978556124878"
  9  = "0.1302808926112106"
  10 = "-0.1944981035472806"
  11 = "-0.3817857436446591"
  12 = "-2.375649628613696E-07"
  13 = "-0.0555296279974082"
  14 = "3.829984367986761E-07"
  15 = "-1.604754923945073E-07"
  16 = "0.009398958989038461"
  17 = "-0.07651818316594991"
  18 = "2.770877186031306E-07"
  19 = "0.2010531357750048"
  20 = "-0.2003621554241067"
}

$firstDataCol = 2   # column B
$lastCol      = 11  # column K - right edge of the sheet's data block

for ($r = 2; $r -le 20; $r++) {

    # Read the row's existing values (columns B..K) before overwriting
    # anything, so the shift-right doesn't clobber data it still needs.
    $existing = @()
    for ($c = $firstDataCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -eq $null) { break }
        $existing += $cell.Value2
    }

    # Insert the corrected naive-forecaster value at column B.
    $ws.Cells.Item($r, $firstDataCol).Value = [double]$newB[$r]

    # Shift the previously-read values one column to the right,
    # starting at column C; anything that would land past column K
    # falls off the sheet's right edge (dropped), same as Excel's
    # own "Insert Cells, Shift Right" behaviour on a bounded range.
    for ($i = 0; $i -lt $existing.Count; $i++) {
        $destCol = $firstDataCol + 1 + $i
        if ($destCol -gt $lastCol) { break }
        $ws.Cells.Item($r, $destCol).Value = $existing[$i]
    }
}
